# Update login credentials in the loginData sheet and trim the
# worksheet back down to a single credential row (rows 2-3 share the
# same username/password), removing the now-unused extra rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loginData")

# Rename the existing credential pair (row 2) to the new values.
$ws.Range("A2").Value = "mngr353180"
$ws.Range("B2").Value = "nerynYt"

# Row 3 repeats the same (renamed) credential pair.
$ws.Range("A3").Value = "mngr353180"
$ws.Range("B3").Value = "nerynYt"

# Remove the now-unused rows 4-6 entirely, shrinking the sheet to A1:B3.
$ws.Rows("4:6").Delete()

# Update the active selection to match the trimmed sheet.
$ws.Range("A3").Select()
